$wb = $excel.ActiveWorkbook

# Insert a new worksheet "6 classes" right after "Accuracy" (before "RandomForest")
$wsAccuracy = $wb.Worksheets.Item("Accuracy")
$ws = $wb.Worksheets.Add($null, $wsAccuracy)
$ws.Name = "6 classes"

$ws.Range("A1").Value = 55.956200000000003
$ws.Range("A2").Value = 56.600099999999998
$ws.Range("A3").Value = 58.8324
$ws.Range("A4").Value = 58.188499999999998
$ws.Range("A5").Value = 61.622700000000002
$ws.Range('A7').Value = 'Classifier Model'
$ws.Range('A8').Value = 'J48 pruned tree'
$ws.Range('A9').Value = '------------------'
$ws.Range('A11').Value = 'COUNT(blunder) = ''(-inf-0.5]'''
$ws.Range('A12').Value = '|   COUNT(misses) = ''(-inf-0.5]'''
$ws.Range('A13').Value = '|   |   COUNT(mistake) = ''(-inf-0.5]'''
$ws.Range('A14').Value = '|   |   |   COUNT(error) = ''(-inf-0.5]'''
$ws.Range('A15').Value = '|   |   |   |   COUNT(interesting) = ''(-inf-0.5]'''
$ws.Range('A16').Value = '|   |   |   |   |   COUNT(loses) = ''(-inf-0.5]'''
$ws.Range('A17').Value = '|   |   |   |   |   |   COUNT(too) = ''(-inf-0.5]'''
$ws.Range('A18').Value = '|   |   |   |   |   |   |   COUNT(should) = ''(-inf-0.5]'''
$ws.Range('A19').Value = '|   |   |   |   |   |   |   |   COUNT(better) = ''(-inf-0.5]'''
$ws.Range('A20').Value = '|   |   |   |   |   |   |   |   |   COUNT(refuted) = ''(-inf-0.5]'''
$ws.Range('A21').Value = '|   |   |   |   |   |   |   |   |   |   COUNT(consideration) = ''(-inf-0.5]'''
$ws.Range('A22').Value = '|   |   |   |   |   |   |   |   |   |   |   COUNT(_blank__carriage_return__new_line_) = ''(-inf-1.5]'': 2 (4060.0/1576.0)'
$ws.Range('A23').Value = '|   |   |   |   |   |   |   |   |   |   |   COUNT(_blank__carriage_return__new_line_) = ''(1.5-inf)'''
$ws.Range('A24').Value = '|   |   |   |   |   |   |   |   |   |   |   |   COUNT(was) = ''(-inf-0.5]'''
$ws.Range('A25').Value = '|   |   |   |   |   |   |   |   |   |   |   |   |   COUNT(.) = ''(-inf-0.5]'': 2 (0.0)'
$ws.Range('A26').Value = '|   |   |   |   |   |   |   |   |   |   |   |   |   COUNT(.) = ''(0.5-3.5]'': 2 (11.0/5.0)'
$ws.Range('A27').Value = '|   |   |   |   |   |   |   |   |   |   |   |   |   COUNT(.) = ''(3.5-inf)'': 5 (17.0/10.0)'
$ws.Range('A28').Value = '|   |   |   |   |   |   |   |   |   |   |   |   COUNT(was) = ''(0.5-inf)'': 4 (8.0/4.0)'
$ws.Range('A29').Value = '|   |   |   |   |   |   |   |   |   |   COUNT(consideration) = ''(0.5-inf)'': 3 (9.0)'
$ws.Range('A30').Value = '|   |   |   |   |   |   |   |   |   COUNT(refuted) = ''(0.5-inf)'': 5 (14.0/2.0)'
$ws.Range('A31').Value = '|   |   |   |   |   |   |   |   COUNT(better) = ''(0.5-inf)'''
$ws.Range('A32').Value = '|   |   |   |   |   |   |   |   |   COUNT(was) = ''(-inf-0.5]'': 2 (103.0/62.0)'
$ws.Range('A33').Value = '|   |   |   |   |   |   |   |   |   COUNT(was) = ''(0.5-inf)'''
$ws.Range('A34').Value = '|   |   |   |   |   |   |   |   |   |   COUNT(this) = ''(-inf-0.5]'': 4 (24.0/14.0)'
$ws.Range('A35').Value = '|   |   |   |   |   |   |   |   |   |   COUNT(this) = ''(0.5-inf)'': 5 (14.0/8.0)'
$ws.Range('A36').Value = '|   |   |   |   |   |   |   COUNT(should) = ''(0.5-inf)'''
$ws.Range('A37').Value = '|   |   |   |   |   |   |   |   COUNT(sacrifice) = ''(-inf-0.5]'''
$ws.Range('A38').Value = '|   |   |   |   |   |   |   |   |   COUNT(_carriage_return__new_line_) = ''(-inf-1.5]'''
$ws.Range('A39').Value = '|   |   |   |   |   |   |   |   |   |   COUNT(.) = ''(-inf-0.5]'': 3 (19.0/8.0)'
$ws.Range('A40').Value = '|   |   |   |   |   |   |   |   |   |   COUNT(.) = ''(0.5-3.5]'''
$ws.Range('A41').Value = '|   |   |   |   |   |   |   |   |   |   |   COUNT(this) = ''(-inf-0.5]'''
$ws.Range('A42').Value = '|   |   |   |   |   |   |   |   |   |   |   |   COUNT(was) = ''(-inf-0.5]'': 4 (25.0/14.0)'
$ws.Range('A43').Value = '|   |   |   |   |   |   |   |   |   |   |   |   COUNT(was) = ''(0.5-inf)'': 3 (3.0/2.0)'
$ws.Range('A44').Value = '|   |   |   |   |   |   |   |   |   |   |   COUNT(this) = ''(0.5-inf)'': 5 (10.0/5.0)'
$ws.Range('A45').Value = '|   |   |   |   |   |   |   |   |   |   COUNT(.) = ''(3.5-inf)'''
$ws.Range('A46').Value = '|   |   |   |   |   |   |   |   |   |   |   COUNT(better) = ''(-inf-0.5]'': 5 (4.0/2.0)'
$ws.Range('A47').Value = '|   |   |   |   |   |   |   |   |   |   |   COUNT(better) = ''(0.5-inf)'': 2 (2.0/1.0)'
$ws.Range('A48').Value = '|   |   |   |   |   |   |   |   |   COUNT(_carriage_return__new_line_) = ''(1.5-2.5]'': 2 (6.0/3.0)'
$ws.Range('A49').Value = '|   |   |   |   |   |   |   |   |   COUNT(_carriage_return__new_line_) = ''(2.5-inf)'''
$ws.Range('A50').Value = '|   |   |   |   |   |   |   |   |   |   COUNT(.) = ''(-inf-0.5]'': 2 (0.0)'
$ws.Range('A51').Value = '|   |   |   |   |   |   |   |   |   |   COUNT(.) = ''(0.5-3.5]'': 2 (9.0/3.0)'
$ws.Range('A52').Value = '|   |   |   |   |   |   |   |   |   |   COUNT(.) = ''(3.5-inf)'': 4 (4.0/1.0)'
$ws.Range('A53').Value = '|   |   |   |   |   |   |   |   COUNT(sacrifice) = ''(0.5-inf)'': 1 (2.0/1.0)'
$ws.Range('A54').Value = '|   |   |   |   |   |   COUNT(too) = ''(0.5-inf)'': 5 (85.0/49.0)'
$ws.Range('A55').Value = '|   |   |   |   |   COUNT(loses) = ''(0.5-inf)'''
$ws.Range('A56').Value = '|   |   |   |   |   |   COUNT(.) = ''(-inf-0.5]'': 6 (11.0/5.0)'
$ws.Range('A57').Value = '|   |   |   |   |   |   COUNT(.) = ''(0.5-3.5]'': 5 (30.0/15.0)'
$ws.Range('A58').Value = '|   |   |   |   |   |   COUNT(.) = ''(3.5-inf)'''
$ws.Range('A59').Value = '|   |   |   |   |   |   |   COUNT(this) = ''(-inf-0.5]'': 2 (4.0/2.0)'
$ws.Range('A60').Value = '|   |   |   |   |   |   |   COUNT(this) = ''(0.5-inf)'': 3 (3.0/1.0)'
$ws.Range('A61').Value = '|   |   |   |   COUNT(interesting) = ''(0.5-inf)'''
$ws.Range('A62').Value = '|   |   |   |   |   COUNT(_blank__carriage_return__new_line_) = ''(-inf-1.5]'': 3 (34.0/11.0)'
$ws.Range('A63').Value = '|   |   |   |   |   COUNT(_blank__carriage_return__new_line_) = ''(1.5-inf)'': 4 (2.0/1.0)'
$ws.Range('A64').Value = '|   |   |   COUNT(error) = ''(0.5-inf)'': 5 (31.0/9.0)'
$ws.Range('A65').Value = '|   |   COUNT(mistake) = ''(0.5-inf)'''
$ws.Range('A66').Value = '|   |   |   COUNT(better) = ''(-inf-0.5]'': 5 (42.0/10.0)'
$ws.Range('A67').Value = '|   |   |   COUNT(better) = ''(0.5-inf)'''
$ws.Range('A68').Value = '|   |   |   |   COUNT(was) = ''(-inf-0.5]'': 2 (3.0/1.0)'
$ws.Range('A69').Value = '|   |   |   |   COUNT(was) = ''(0.5-inf)'': 6 (3.0/1.0)'
$ws.Range('A70').Value = '|   COUNT(misses) = ''(0.5-inf)'': 5 (27.0/12.0)'
$ws.Range('A71').Value = 'COUNT(blunder) = ''(0.5-inf)'''
$ws.Range('A72').Value = '|   COUNT(game-losing) = ''(-inf-0.5]'': 5 (35.0/17.0)'
$ws.Range('A73').Value = '|   COUNT(game-losing) = ''(0.5-inf)'': 6 (5.0)'
$ws.Range('A75').Value = 'Number of Leaves  : '
$ws.Range("B75").Value = 35
$ws.Range('A77').Value = 'Size of the tree : '
$ws.Range("B77").Value = 64
$ws.Range('A80').Value = 'Time taken to build model: 28.05 seconds'
$ws.Range('A82').Value = '''=== Stratified cross-validation ==='
$ws.Range('A83').Value = '''=== Summary ==='
$ws.Range('A85').Value = 'Correctly Classified Instances        2711               58.1885 %'
$ws.Range('A86').Value = 'Incorrectly Classified Instances      1948               41.8115 %'
$ws.Range('A87').Value = 'Kappa statistic                          0.1314'
$ws.Range('A88').Value = 'Mean absolute error                      0.2003'
$ws.Range('A89').Value = 'Root mean squared error                  0.3197'
$ws.Range('A90').Value = 'Relative absolute error                 94.2179 %'
$ws.Range('A91').Value = 'Root relative squared error             98.0835 %'
$ws.Range('A92').Value = 'Total Number of Instances             4659     '
$ws.Range('A94').Value = '''=== Detailed Accuracy By Class ==='
$ws.Range('A96').Value = '                 TP Rate  FP Rate  Precision  Recall   F-Measure  MCC      ROC Area  PRC Area  Class'
$ws.Range('A97').Value = '                 0,000    0,001    0,000      0,000    0,000      -0,006   0,529     0,038     1'
$ws.Range('A98').Value = '                 0,965    0,827    0,597      0,965    0,738      0,233    0,575     0,598     2'
$ws.Range('A99').Value = '                 0,068    0,011    0,483      0,068    0,120      0,143    0,546     0,180     3'
$ws.Range('A100').Value = '                 0,050    0,011    0,303      0,050    0,086      0,093    0,558     0,125     4'
$ws.Range('A101').Value = '                 0,183    0,036    0,467      0,183    0,263      0,223    0,592     0,276     5'
$ws.Range('A102').Value = '                 0,033    0,002    0,353      0,033    0,061      0,099    0,540     0,087     6'
$ws.Range('A103').Value = 'Weighted Avg.    0,582    0,471    0,507      0,582    0,477      0,194    0,569     0,415     '
$ws.Range('A105').Value = '''=== Confusion Matrix ==='
$ws.Range('A107').Value = '    a    b    c    d    e    f   <-- classified as'
$ws.Range('A108').Value = '    0  162    0    2    4    0 |    a = 1'
$ws.Range('A109').Value = '    1 2517   23   15   48    3 |    b = 2'
$ws.Range('A110').Value = '    1  551   42    5   16    0 |    c = 3'
$ws.Range('A111').Value = '    0  326    7   20   47    0 |    d = 4'
$ws.Range('A112').Value = '    2  522   11   20  126    8 |    e = 5'
$ws.Range('A113').Value = '    0  137    4    4   29    6 |    f = 6'
$ws.Range('A116').Value = '''=== Run information ==='
$ws.Range('A118').Value = 'Scheme:       weka.classifiers.trees.RandomForest -P 100 -I 100 -num-slots 1 -K 0 -M 1.0 -V 0.001 -S 1'
$ws.Range('A119').Value = 'Relation:     comment'
$ws.Range('A120').Value = 'Instances:    4659'
$ws.Range('A121').Value = 'Attributes:   6365'
$ws.Range('A122').Value = '              [list of attributes omitted]'
$ws.Range('A123').Value = 'Test mode:    10-fold cross-validation'
$ws.Range('A125').Value = '''=== Classifier model (full training set) ==='
$ws.Range('A127').Value = 'RandomForest'
$ws.Range('A129').Value = 'Bagging with 100 iterations and base learner'
$ws.Range('A131').Value = 'weka.classifiers.trees.RandomTree -K 0 -M 1.0 -V 0.001 -S 1 -do-not-check-capabilities'
$ws.Range('A133').Value = 'Time taken to build model: 118.92 seconds'
$ws.Range('A135').Value = '''=== Stratified cross-validation ==='
$ws.Range('A136').Value = '''=== Summary ==='
$ws.Range('A138').Value = 'Correctly Classified Instances        2871               61.6227 %'
$ws.Range('A139').Value = 'Incorrectly Classified Instances      1788               38.3773 %'
$ws.Range('A140').Value = 'Kappa statistic                          0.2049'
$ws.Range('A141').Value = 'Mean absolute error                      0.1791'
$ws.Range('A142').Value = 'Root mean squared error                  0.2968'
$ws.Range('A143').Value = 'Relative absolute error                 84.2655 %'
$ws.Range('A144').Value = 'Root relative squared error             91.0506 %'
$ws.Range('A145').Value = 'Total Number of Instances             4659     '
$ws.Range('A147').Value = '''=== Detailed Accuracy By Class ==='
$ws.Range('A149').Value = '                 TP Rate  FP Rate  Precision  Recall   F-Measure  MCC      ROC Area  PRC Area  Class'
$ws.Range('A150').Value = '                 0,000    0,001    0,000      0,000    0,000      -0,007   0,690     0,101     1'
$ws.Range('A151').Value = '                 0,982    0,805    0,608      0,982    0,751      0,299    0,776     0,795     2'
$ws.Range('A152').Value = '                 0,172    0,010    0,726      0,172    0,279      0,316    0,761     0,419     3'
$ws.Range('A153').Value = '                 0,083    0,005    0,589      0,083    0,145      0,198    0,725     0,259     4'
$ws.Range('A154').Value = '                 0,179    0,016    0,661      0,179    0,281      0,295    0,787     0,458     5'
$ws.Range('A155').Value = '                 0,272    0,001    0,907      0,272    0,419      0,488    0,821     0,418     6'
$ws.Range('A156').Value = 'Weighted Avg.    0,616    0,454    0,619      0,616    0,527      0,288    0,770     0,610     '
$ws.Range('A158').Value = '''=== Confusion Matrix ==='
$ws.Range('A160').Value = '    a    b    c    d    e    f   <-- classified as'
$ws.Range('A161').Value = '    0  167    1    0    0    0 |    a = 1'
$ws.Range('A162').Value = '    6 2560   28    4    9    0 |    b = 2'
$ws.Range('A163').Value = '    0  506  106    2    1    0 |    c = 3'
$ws.Range('A164').Value = '    0  332    5   33   30    0 |    d = 4'
$ws.Range('A165').Value = '    0  539    5   17  123    5 |    e = 5'
$ws.Range('A166').Value = '    0  107    1    0   23   49 |    f = 6'

$ws.Activate()
$ws.Range("A116:A167").Select()
$ws.Range("A116").Activate()
